# Trade #68 closed at 2026-02-16 21:35:18 - momentum DOWN +0.000%
#
# This records two things that happened in the live-trading log:
#   1) leadlag Trade #48 (opened 21:30:18, row 38 on the "leadlag" sheet)
#      finished closing out - exit price, P&L and duration land in the
#      "leadlag" sheet and get appended as a new row on "All Trades".
#   2) momentum Trade #68 (21:35:18) opens - appended as a new row on the
#      "momentum" sheet.
#   3) The rollup numbers on "Summary" and "Comparison" shift accordingly.
#
# NOTE: many of the text columns in this workbook hold strings that LOOK
# like numbers/dates/percentages (e.g. "64.6%", "2.43", "2026-02-16").
# Assigning those bare to Range.Value lets Excel's type-inference turn them
# into real numbers/dates, which is not what the source data has (they are
# plain text cells). Prefixing with a leading apostrophe forces text entry,
# exactly like typing ' into the formula bar - the apostrophe itself is not
# stored as part of the value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - update aggregate OVERALL and leadlag rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C2").Value = 48
$summary.Range("D2").Value = "'64.6%"
$summary.Range("D2").Style = "Normal"
$summary.Range("E2").Value = "'+11.3229%"
$summary.Range("E2").Style = "Normal"
$summary.Range("F2").Value = "'+0.2359%"
$summary.Range("F2").Style = "Normal"
$summary.Range("E3").Value = "'+7.1213%"
$summary.Range("E3").Style = "Normal"
$summary.Range("F3").Value = "'+0.1344%"
$summary.Range("F3").Style = "Normal"

# ---------------------------------------------------------------------
# 2) leadlag sheet - trade #48 (row 38) closes out
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")
$leadlag.Range("G38").Value = 68480.61685200001
$leadlag.Range("H38").Value = "CLOSED"
$leadlag.Range("I38").Value = -0.3451
$leadlag.Range("J38").Value = -3.45
$leadlag.Range("M38").Value = "time_exit_5min"
$leadlag.Range("N38").Value = 5

# ---------------------------------------------------------------------
# 3) momentum sheet - new trade #68 opened (row 16)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A16").Value = 68
$momentum.Range("B16").Value = "'2026-02-16"
$momentum.Range("B16").Style = "Normal"
$momentum.Range("C16").Value = "'21:35:18"
$momentum.Range("C16").Style = "Normal"
$momentum.Range("D16").Value = "momentum"
$momentum.Range("E16").Value = "DOWN"
$momentum.Range("F16").Value = 68557.14999999999
$momentum.Range("G16").Value = "'"
$momentum.Range("G16").Style = "Normal"
$momentum.Range("H16").Value = "OPEN"
$momentum.Range("I16").Value = 0
$momentum.Range("J16").Value = 0
$momentum.Range("K16").Value = 0.9
$momentum.Range("L16").Value = "Downward momentum: -0.223% over 10 samples"
$momentum.Range("M16").Value = "'"
$momentum.Range("M16").Style = "Normal"
$momentum.Range("N16").Value = 0

# ---------------------------------------------------------------------
# 4) All Trades sheet - append trade #48 (leadlag) as CLOSED (row 49)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A49").Value = 48
$allTrades.Range("B49").Value = "'2026-02-16"
$allTrades.Range("B49").Style = "Normal"
$allTrades.Range("C49").Value = "'21:30:18"
$allTrades.Range("C49").Style = "Normal"
$allTrades.Range("D49").Value = "leadlag"
$allTrades.Range("E49").Value = "UP"
$allTrades.Range("F49").Value = 68717.73
$allTrades.Range("G49").Value = 68480.61685200001
$allTrades.Range("H49").Value = "CLOSED"
$allTrades.Range("I49").Value = -0.3451
$allTrades.Range("J49").Value = -3.45
$allTrades.Range("K49").Value = 0.75
$allTrades.Range("L49").Value = "Binance leading with 0.115% move"
$allTrades.Range("M49").Value = "time_exit_5min"
$allTrades.Range("N49").Value = 5

# ---------------------------------------------------------------------
# 5) Comparison sheet - updated leadlag profit factor / avg loss %
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")
$comparison.Range("D2").Value = "'2.43"
$comparison.Range("D2").Style = "Normal"
$comparison.Range("F2").Value = "'-0.3316%"
$comparison.Range("F2").Style = "Normal"
